$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Metrics" text for the "Using Pirate Metrics..." row to mention retention as well.
$ws.Range("E2").Value = "*pirate metrics (community involvement and retention)"

# Move the view back to the top of the sheet and select E5 (matches author's final cursor position).
$ws.Range("A1").Select()
$ws.Range("E5").Select()
